# Fernandez-Galeano Publication Impact workbook:
# - rename existing sheet to "Journal Metrics"
# - add a new "Article Metrics" sheet after it (and make it the active sheet)
# - add the Journal-Level Metrics header row (TITLE / YEAR / JOURNAL_TITLE / JCR_JIF / SCOPUS_CITESCORE)
# - size up the JOURNAL_TITLE / SCOPUS_CITESCORE columns so the headers aren't clipped

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Journal Metrics"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Article Metrics"

$ws1.Range("A1").Value = "TITLE"
$ws1.Range("B1").Value = "YEAR"
$ws1.Range("C1").Value = "JOURNAL_TITLE"
$ws1.Range("D1").Value = "JCR_JIF"
$ws1.Range("E1").Value = "SCOPUS_CITESCORE"

# Target stored widths are ~23.66 ("JOURNAL_TITLE") and ~22.66 ("SCOPUS_CITESCORE")
# characters; ColumnWidth assignments land on the nearest representable width.
$ws1.Columns.Item(3).ColumnWidth = 22.8
$ws1.Columns.Item(5).ColumnWidth = 21.8

$ws1.Range("D1").Select() | Out-Null

$ws2.Range("D32").Select() | Out-Null
$ws2.Activate() | Out-Null
